$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "Amphenol-SAA"
$ws.Range("D14").Value = "PIOV008NRAA-100"
$ws.Range("E14").Value = "PIOV008NRAA-100"
$ws.Range("F14").Value = "Mouser"
$ws.Range("G14").Value = "523-PIOV008NRAA-100"
$ws.Range("H14").Value = 2.69
$ws.Range("J14").Value = "https://www.mouser.de/ProductDetail/Amphenol-SAA/PIOV008NRAA-100?qs=GedFDFLaBXFCaCiGvxFhnA%3D%3D"
